$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45400 to 45402 for rows 2-29
$ws.Range("C2:C29").Value = 45402
